# Auto-generated script applying the Durandal_Profits.xlsx diff
# Updates per-row numeric values (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 2697
$ws.Range("I51").Value = 1799.5
$ws.Range("K51").Value = 1799.5
$ws.Range("M51").Value = -1315.5

# Row 64
$ws.Range("H64").Value = 2977.547
$ws.Range("I64").Value = 2973.4883
$ws.Range("J64").Value = 2995
$ws.Range("K64").Value = 2973.4883
$ws.Range("L64").Value = 2995
$ws.Range("M64").Value = -2725.4883
$ws.Range("N64").Value = -3491

# Row 67
$ws.Range("H67").Value = 2977.547
$ws.Range("I67").Value = 2973.4883
$ws.Range("J67").Value = 2995
$ws.Range("K67").Value = 2973.4883
$ws.Range("L67").Value = 2995
$ws.Range("M67").Value = -2115.4883
$ws.Range("N67").Value = -4711

# Row 100
$ws.Range("H100").Value = 17546306
$ws.Range("I100").Value = 23811234
$ws.Range("J100").Value = 4509.2
$ws.Range("K100").Value = 23811234
$ws.Range("L100").Value = 4509.2
$ws.Range("M100").Value = -23810693
$ws.Range("N100").Value = -5591.2

# Row 135
$ws.Range("H135").Value = 2174.2292
$ws.Range("I135").Value = 984.48834
$ws.Range("J135").Value = 12406
$ws.Range("K135").Value = 8860.395060000001
$ws.Range("L135").Value = 111654
$ws.Range("M135").Value = -6325.395060000001
$ws.Range("N135").Value = -116724

# Row 137
$ws.Range("H137").Value = 886.6486
$ws.Range("I137").Value = 765.5909
$ws.Range("J137").Value = 1064.2
$ws.Range("K137").Value = 2296.7727
$ws.Range("L137").Value = 3192.6
$ws.Range("M137").Value = 253.2273
$ws.Range("N137").Value = -8292.6

# Row 138
$ws.Range("H138").Value = 3969.3518
$ws.Range("I138").Value = 2423.5
$ws.Range("J138").Value = 5634.115
$ws.Range("K138").Value = 7270.5
$ws.Range("L138").Value = 16902.345
$ws.Range("M138").Value = -2130.5
$ws.Range("N138").Value = -27182.345

# Row 141
$ws.Range("H141").Value = 3483.8635
$ws.Range("I141").Value = 3554.7368
$ws.Range("J141").Value = 3035
$ws.Range("K141").Value = 10664.2104
$ws.Range("L141").Value = 9105
$ws.Range("M141").Value = -5484.2104
$ws.Range("N141").Value = -19465

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3530.5454
$ws.Range("I2").Value = 1852.875
$ws.Range("J2").Value = 8004.3335
$ws.Range("K2").Value = 1852.875
$ws.Range("L2").Value = 8004.3335
$ws.Range("M2").Value = -1739.875
$ws.Range("N2").Value = -8230.333500000001

# Row 32
$ws.Range("H32").Value = 360332.62
$ws.Range("I32").Value = 2543.1487
$ws.Range("J32").Value = 4773069.5
$ws.Range("K32").Value = 2543.1487
$ws.Range("L32").Value = 4773069.5
$ws.Range("M32").Value = -2256.1487
$ws.Range("N32").Value = -4773643.5

# Row 81
$ws.Range("H81").Value = 20172.5
$ws.Range("J81").Value = 181
$ws.Range("L81").Value = 181
$ws.Range("N81").Value = -2177

# Row 84
$ws.Range("H84").Value = 20172.5
$ws.Range("J84").Value = 181
$ws.Range("L84").Value = 543
$ws.Range("N84").Value = -10527

# Row 97
$ws.Range("H97").Value = 1409.7188
$ws.Range("I97").Value = 868.2727
$ws.Range("J97").Value = 2600.9
$ws.Range("K97").Value = 868.2727
$ws.Range("L97").Value = 2600.9
$ws.Range("M97").Value = -372.2727
$ws.Range("N97").Value = -3592.9

# Row 116
$ws.Range("H116").Value = 3530.5454
$ws.Range("I116").Value = 1852.875
$ws.Range("J116").Value = 8004.3335
$ws.Range("K116").Value = 1852.875
$ws.Range("L116").Value = 8004.3335
$ws.Range("M116").Value = 441.125
$ws.Range("N116").Value = -12592.3335

# Row 135
$ws.Range("H135").Value = 58000
$ws.Range("J135").Value = 58000
$ws.Range("L135").Value = 58000
$ws.Range("N135").Value = -68140

# Row 137
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

# Row 138
$ws.Range("H138").Value = 23760
$ws.Range("J138").Value = 23760
$ws.Range("L138").Value = 23760
$ws.Range("N138").Value = -34040

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3530.5454
$ws.Range("I3").Value = 1852.875
$ws.Range("J3").Value = 8004.3335
$ws.Range("K3").Value = 1852.875
$ws.Range("L3").Value = 8004.3335
$ws.Range("M3").Value = -1738.875
$ws.Range("N3").Value = -8232.333500000001

# Row 134
$ws.Range("H134").Value = 4073.5789
$ws.Range("I134").Value = 1307.6471
$ws.Range("J134").Value = 8162.3477
$ws.Range("K134").Value = 3922.9413
$ws.Range("L134").Value = 24487.0431
$ws.Range("M134").Value = -1387.9413
$ws.Range("N134").Value = -29557.0431

# Row 137
$ws.Range("H137").Value = 37500
$ws.Range("J137").Value = 37500
$ws.Range("L137").Value = 37500
$ws.Range("N137").Value = -47700

# Row 138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# Row 140
$ws.Range("H140").Value = 61900
$ws.Range("J140").Value = 61900
$ws.Range("L140").Value = 61900
$ws.Range("N140").Value = -72260

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9261438
$ws.Range("I31").Value = 13159743
$ws.Range("J31").Value = 2964.4375
$ws.Range("K31").Value = 13159743
$ws.Range("L31").Value = 2964.4375
$ws.Range("M31").Value = -13159448
$ws.Range("N31").Value = -3554.4375

# Row 34
$ws.Range("H34").Value = 9261438
$ws.Range("I34").Value = 13159743
$ws.Range("J34").Value = 2964.4375
$ws.Range("K34").Value = 13159743
$ws.Range("L34").Value = 2964.4375
$ws.Range("M34").Value = -13159541
$ws.Range("N34").Value = -3368.4375

# Row 132
$ws.Range("H132").Value = 1072.723
$ws.Range("I132").Value = 1010.7593
$ws.Range("J132").Value = 1376.909
$ws.Range("K132").Value = 3032.2779
$ws.Range("L132").Value = 4130.727000000001
$ws.Range("M132").Value = -502.2779
$ws.Range("N132").Value = -9190.727000000001

# Row 134
$ws.Range("H134").Value = 1593.1082
$ws.Range("I134").Value = 1564.4
$ws.Range("J134").Value = 1716.1428
$ws.Range("K134").Value = 4693.200000000001
$ws.Range("L134").Value = 5148.428400000001
$ws.Range("M134").Value = -2158.200000000001
$ws.Range("N134").Value = -10218.4284

# Row 135
$ws.Range("H135").Value = 41650
$ws.Range("J135").Value = 41650
$ws.Range("L135").Value = 41650
$ws.Range("N135").Value = -51790

$ws = $wb.Worksheets.Item("CUL")
# Row 114
$ws.Range("H114").Value = 881.3158
$ws.Range("I114").Value = 286.63635
$ws.Range("J114").Value = 1699
$ws.Range("K114").Value = 859.90905
$ws.Range("L114").Value = 5097
$ws.Range("M114").Value = 2394.09095
$ws.Range("N114").Value = -11605

# Row 131
$ws.Range("H131").Value = 7693162
$ws.Range("I131").Value = 834.1667
$ws.Range("J131").Value = 9434821
$ws.Range("K131").Value = 2502.5001
$ws.Range("L131").Value = 28304463
$ws.Range("M131").Value = 2537.4999
$ws.Range("N131").Value = -28314543

# Row 132
$ws.Range("H132").Value = 1301.75
$ws.Range("J132").Value = 1379.8823
$ws.Range("L132").Value = 12418.9407
$ws.Range("N132").Value = -17478.9407

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 14812416
$ws.Range("I70").Value = 35169430
$ws.Range("J70").Value = 7314.3184
$ws.Range("K70").Value = 35169430
$ws.Range("L70").Value = 7314.3184
$ws.Range("M70").Value = -35169160
$ws.Range("N70").Value = -7854.3184

# Row 73
$ws.Range("H73").Value = 14812416
$ws.Range("I73").Value = 35169430
$ws.Range("J73").Value = 7314.3184
$ws.Range("K73").Value = 35169430
$ws.Range("L73").Value = 7314.3184
$ws.Range("M73").Value = -35168494
$ws.Range("N73").Value = -9186.3184

# Row 80
$ws.Range("H80").Value = 3299.3333
$ws.Range("I80").Value = 2998.75
$ws.Range("J80").Value = 3736.5454
$ws.Range("K80").Value = 2998.75
$ws.Range("L80").Value = 3736.5454
$ws.Range("M80").Value = -2000.75
$ws.Range("N80").Value = -5732.5454

# Row 83
$ws.Range("H83").Value = 3299.3333
$ws.Range("I83").Value = 2998.75
$ws.Range("J83").Value = 3736.5454
$ws.Range("K83").Value = 14993.75
$ws.Range("L83").Value = 18682.727
$ws.Range("M83").Value = -10001.75
$ws.Range("N83").Value = -28666.727

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 1786.4266
$ws.Range("I132").Value = 1751.5614
$ws.Range("K132").Value = 5254.6842
$ws.Range("M132").Value = -2724.6842

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 17362044
$ws.Range("I132").Value = 18940270
$ws.Range("K132").Value = 56820810
$ws.Range("M132").Value = -56818280

# Row 136
$ws.Range("H136").Value = 613.9726000000001
$ws.Range("I136").Value = 546.0645
$ws.Range("J136").Value = 996.7273
$ws.Range("K136").Value = 1638.1935
$ws.Range("L136").Value = 2990.1819
$ws.Range("M136").Value = 911.8065000000001
$ws.Range("N136").Value = -8090.1819

# Row 138
$ws.Range("H138").Value = 94800
$ws.Range("J138").Value = 94800
$ws.Range("L138").Value = 94800
$ws.Range("N138").Value = -105080
